$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '65.342.59'
$ws.Range("E2").Value = '  -5.98%  '
$ws.Range("D3").Value = '3.523.96'
$ws.Range("E3").Value = '  -5.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -8.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("D7").Value = '3.516.78'
$ws.Range("E7").Value = '  -5.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.85%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.660'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.51'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.142'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -13.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000252'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -15.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -10.56%  '
$ws.Range("D15").Value = '4.094.18'
$ws.Range("D16").Value = '3.531.17'
$ws.Range("E16").Value = '  -5.66%  '
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -7.55%  '
$ws.Range("D19").Value = '65.318.42'
$ws.Range("E19").Value = '  -5.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.86%  '
$ws.Range("E21").Value = '  -8.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.09'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -11.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.60%  '
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("E28").Value = '  -7.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.68%  '
$ws.Range("E32").Value = '  -10.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '605.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '64.59'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.37%  '
$ws.Range("E36").Value = '  -7.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.75%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").Value = '0.0₃0739'
$ws.Range("E40").Value = '  -17.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.368'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -10.34%  '
$ws.Range("E42").Value = '  -7.62%  '
$ws.Range("D43").Value = '2.843.87'
$ws.Range("E43").Value = '  +0.88%  '
$ws.Range("E44").Value = '  -11.30%  '
$ws.Range("E45").Value = '  -9.17%  '
$ws.Range("E46").Value = '  -10.42%  '
$ws.Range("E47").Value = '  -5.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '137.25'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -10.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.16'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -12.23%  '
